# Templates_Feature.xlsx edit:
#  - Add a new "URL" worksheet at the end of the workbook with a
#    "Parcel URL" label and an actual hyperlink to the UAT login page.
#  - Remove the stray row 4 ("1567") from "Template Setup" and move the
#    selection there.

$wb = $excel.ActiveWorkbook

# --- Template Setup: drop row 4, move selection -------------------------
$setupWs = $wb.Worksheets.Item("Template Setup")
$setupWs.Rows.Item(4).Delete()
$setupWs.Range("D10").Select()

# --- New "URL" worksheet, inserted after the last existing sheet --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$urlWs = $wb.Worksheets.Add($null, $lastSheet)
$urlWs.Name = "URL"

$urlWs.Columns.Item(1).ColumnWidth = 51

$urlWs.Range("A1").Value = "Parcel URL"
$urlWs.Range("A1").Font.Bold = $true
$urlWs.Range("A1").HorizontalAlignment = -4108
$urlWs.Range("A1").VerticalAlignment = -4108

$urlWs.Range("A2").Value = "https://uat.parcelplatform.com/reporting/login.php"
$urlWs.Hyperlinks.Add($urlWs.Range("A2"), "https://uat.parcelplatform.com/reporting/login.php")

$urlWs.Range("D13").Select()
